# Macroferia Regional de Talca - Berenjena: add one new weekly record.
# Insert a brand-new row at sheet row 82, pushing the existing rows 82-140
# down to 83-141, then populate the new row with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 82..140 down to 83..141 by inserting a fresh row at 82.
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with this week's record.
$ws.Cells.Item(82, 1).Value = 5
$ws.Cells.Item(82, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(82, 3).Value = "Maule"
$ws.Cells.Item(82, 4).Value = 44978
$ws.Cells.Item(82, 5).Value = 7
$ws.Cells.Item(82, 6).Value = 100112001
$ws.Cells.Item(82, 7).Value = "Berenjena"
$ws.Cells.Item(82, 8).Value = "Sin especificar"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 150
$ws.Cells.Item(82, 11).Value = 7000
$ws.Cells.Item(82, 12).Value = 7000
$ws.Cells.Item(82, 13).Value = 7000
$ws.Cells.Item(82, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(82, 15).Value = "Región del Maule"
$ws.Cells.Item(82, 16).Value = 140
$ws.Cells.Item(82, 17).Value = 50
$ws.Cells.Item(82, 18).Value = "Hortaliza"
